$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41; existing rows 41..76 shift down to 42..77
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new price-report record
$ws.Cells.Item(41, 1).Value = 5
$ws.Cells.Item(41, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(41, 3).Value = "Maule"
$ws.Cells.Item(41, 4).Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44589)
$ws.Cells.Item(41, 5).Value = 7
$ws.Cells.Item(41, 6).Value = "Fruta"
$ws.Cells.Item(41, 7).Value = 100103
$ws.Cells.Item(41, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(41, 9).Value = 100103002
$ws.Cells.Item(41, 10).Value = "Ciruela"
$ws.Cells.Item(41, 11).Value = "Black Amber"
$ws.Cells.Item(41, 12).Value = "Primera"
$ws.Cells.Item(41, 13).Value = 200
$ws.Cells.Item(41, 14).Value = 9000
$ws.Cells.Item(41, 15).Value = 9000
$ws.Cells.Item(41, 16).Value = 9000
$ws.Cells.Item(41, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(41, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(41, 19).Value = 500
$ws.Cells.Item(41, 20).Value = 18
